$wb = $excel.ActiveWorkbook

# --- Sheet 1: Model Accuracy (-0.65, 0.65, 0.65) ---
$ws1 = $wb.Worksheets.Item("Model Accuracy (-0.65, 0.65, 0.65)")

$ws1.Range("B1").Copy()
$ws1.Range("C1:G1").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("C1").Value = "Market threshold"
$ws1.Range("D1").Value = "Market min"
$ws1.Range("E1").Value = "Market max"
$ws1.Range("F1").Value = "Recall"
$ws1.Range("G1").Value = "Precision"

$ws1.Range("B2").Value = 62.71393643031785
$ws1.Range("C2").Value = 0.05450546436368681
$ws1.Range("D2").Value = -15.55441
$ws1.Range("E2").Value = 15.06418
$ws1.Range("F2").Value = 0
$ws1.Range("G2").Value = 0

$ws1.Range("B3").Value = 37.28606356968216
$ws1.Range("C3").Value = 0.009583939973006913
$ws1.Range("D3").Value = -19.35264
$ws1.Range("E3").Value = 13.70093
$ws1.Range("F3").Value = 2.412868632707775
$ws1.Range("G3").Value = 24.32432432432433

$ws1.Range("B4").Value = 92.29828850855746
$ws1.Range("C4").Value = 0.04158117063764853
$ws1.Range("D4").Value = -18.75314
$ws1.Range("E4").Value = 23.33066
$ws1.Range("F4").Value = 0
$ws1.Range("G4").Value = 0

$ws1.Range("B5").Value = 81.2958435207824
$ws1.Range("C5").Value = 0.02983403801513819
$ws1.Range("D5").Value = -12.78028
$ws1.Range("E5").Value = 12.42348
$ws1.Range("F5").Value = 0
$ws1.Range("G5").Value = 0

$ws1.Range("B6").Value = 94.98777506112469
$ws1.Range("C6").Value = 0.08368817696170747
$ws1.Range("D6").Value = -16.47904
$ws1.Range("E6").Value = 14.94325
$ws1.Range("F6").Value = 0
$ws1.Range("G6").Value = 0

# --- Sheet 2: Confusion Matrix TOTALENERGIES SE ---
$ws2 = $wb.Worksheets.Item("Confusion Matrix TOTALENERGIES SE (-0.65, 0.65, 0.65)")
$ws2.Range("B3").Value = 9
$ws2.Range("C3").Value = 1024
$ws2.Range("D3").Value = 8

# --- Sheet 3: Confusion Matrix FMC CORP ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix FMC CORP (-0.65, 0.65, 0.65)")
$ws3.Range("B2").Value = 9
$ws3.Range("C2").Value = 21
$ws3.Range("D2").Value = 7

$ws3.Range("B3").Value = 321
$ws3.Range("C3").Value = 563
$ws3.Range("D3").Value = 311

$ws3.Range("B4").Value = 43
$ws3.Range("C4").Value = 68
$ws3.Range("D4").Value = 38

# --- Sheet 4: Confusion Matrix BP PLC ---
$ws4 = $wb.Worksheets.Item("Confusion Matrix BP PLC (-0.65, 0.65, 0.65)")
$ws4.Range("B3").Value = 39
$ws4.Range("C3").Value = 1510
$ws4.Range("D3").Value = 42

$ws4.Range("B4").Value = 1
$ws4.Range("C4").Value = 7

# --- Sheet 5: Confusion Matrix STORA ENSO ---
$ws5 = $wb.Worksheets.Item("Confusion Matrix STORA ENSO (-0.65, 0.65, 0.65)")
$ws5.Range("B3").Value = 108
$ws5.Range("C3").Value = 1330
$ws5.Range("D3").Value = 107

$ws5.Range("B4").Value = 2
$ws5.Range("C4").Value = 20

# --- Sheet 6: Confusion Matrix BHP GROUP ---
$ws6 = $wb.Worksheets.Item("Confusion Matrix BHP GROUP (-0.65, 0.65, 0.65)")
$ws6.Range("B3").Value = 4
$ws6.Range("C3").Value = 1554
$ws6.Range("D3").Value = 3
